# "break out stock.yaml completed"
#
# On the "10per change" sheet:
#   1) bsecode column (D11:D21) was being written as literal text
#      ("541729", ...) by the old scraper; re-type those 11 cells as real
#      numbers (values themselves are unchanged).
#   2) A newer chartink run (06/06/2024 13:16:38) was appended as rows
#      22-32, mirroring the nsecode/name/bsecode of rows 11-21 but with
#      refreshed per_chg / close / volume figures. Those bsecode values
#      stay text, same as the rest of the sheet's history.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "10per change"

# --- 1) D11:D21: text "541729" etc. -> numeric 541729 ---------------------
$bsecodes = @(541729, 500510, 512599, 541450, 532921, 500112, 533096, 533278, 532898, 532134, 532155)
for ($i = 0; $i -lt $bsecodes.Length; $i++) {
    $row = 11 + $i
    $ws.Cells.Item($row, 4).Value = $bsecodes[$i]
}

# --- 2) Append the new snapshot, rows 22-32 --------------------------------
# columns: sr, nsecode, name, bsecode, per_chg, close, volume, Date Time
$newRows = @(
    @(1,  "HDFCAMC",    "HDFC Asset Management Company Ltd",             "541729", 1.49,  3760.5,  451506,    "06/06/2024 13:16:38"),
    @(2,  "LT",         "Larsen & Toubro Limited",                       "500510", 2.16,  3482.55, 8374602,   "06/06/2024 13:16:38"),
    @(3,  "ADANIENT",   "Adani Enterprises Limited",                     "512599", 2.26,  3185.65, 5873271,   "06/06/2024 13:16:38"),
    @(4,  "ADANIGREEN", "Adani Green Energy Ltd",                        "541450", 1.99,  1865.2,  1909863,   "06/06/2024 13:16:38"),
    @(5,  "ADANIPORTS", "Adani Ports And Special Economic Zone Limited", "532921", -0.12, 1352.95, 12153840,  "06/06/2024 13:16:38"),
    @(6,  "SBIN",       "State Bank Of India",                           "500112", 3.44,  816.95,  40555649,  "06/06/2024 13:16:38"),
    @(7,  "ADANIPOWER", "Adani Power Limited",                           "533096", 3.12,  749.35,  17408648,  "06/06/2024 13:16:38"),
    @(8,  "COALINDIA",  "Coal India Limited",                            "533278", 2.69,  472.95,  21757857,  "06/06/2024 13:16:38"),
    @(9,  "POWERGRID",  "Power Grid Corporation Of India Limited",       "532898", 0.55,  300.45,  36508044,  "06/06/2024 13:16:38"),
    @(10, "BANKBARODA", "Bank Of Baroda",                                "532134", 3.42,  268.9,   31690680,  "06/06/2024 13:16:38"),
    @(11, "GAIL",       "Gail (india) Limited",                          "532155", 6.53,  207.9,   41316996,  "06/06/2024 13:16:38")
)

$startRow = 22
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]

    # Keep bsecode as text (leading apostrophe forces text-entry, same as
    # how the sheet already stores these codes elsewhere), then drop back
    # to the Normal style so no stray number-format sticks to the cell.
    $dcell = $ws.Cells.Item($r, 4)
    $dcell.Value = "'" + $data[3]
    $dcell.Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
}
